$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 589942.5
$ws.Range("I19").Value = 1112562.9
$ws.Range("J19").Value = 1994.625
$ws.Range("K19").Value = 1112562.9
$ws.Range("L19").Value = 1994.625
$ws.Range("M19").Value = -1112387.9
$ws.Range("N19").Value = -2344.625
$ws.Range("H39").Value = 832.61536
$ws.Range("I39").Value = 164.8
$ws.Range("J39").Value = 1250
$ws.Range("K39").Value = 494.4
$ws.Range("L39").Value = 3750
$ws.Range("M39").Value = -198.4
$ws.Range("N39").Value = -4342
$ws.Range("H88").Value = 5987.875
$ws.Range("I88").Value = 967.6667
$ws.Range("J88").Value = 9000
$ws.Range("K88").Value = 967.6667
$ws.Range("L88").Value = 9000
$ws.Range("M88").Value = -561.6667
$ws.Range("N88").Value = -9812
$ws.Range("H91").Value = 5987.875
$ws.Range("I91").Value = 967.6667
$ws.Range("J91").Value = 9000
$ws.Range("K91").Value = 967.6667
$ws.Range("L91").Value = 9000
$ws.Range("M91").Value = 436.3333
$ws.Range("N91").Value = -11808
$ws.Range("H133").Value = 54899
$ws.Range("J133").Value = 54899
$ws.Range("L133").Value = 54899
$ws.Range("N133").Value = -65019
$ws.Range("H137").Value = 32368.516
$ws.Range("I137").Value = 1201.1818
$ws.Range("J137").Value = 47952.184
$ws.Range("K137").Value = 3603.5454
$ws.Range("L137").Value = 143856.552
$ws.Range("M137").Value = -1053.5454
$ws.Range("N137").Value = -148956.552
$ws.Range("H138").Value = 3307.7727
$ws.Range("I138").Value = 4051.077
$ws.Range("J138").Value = 2234.111
$ws.Range("K138").Value = 12153.231
$ws.Range("L138").Value = 6702.333
$ws.Range("M138").Value = -7013.231
$ws.Range("N138").Value = -16982.333

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3111.3333
$ws.Range("I32").Value = 2498.0488
$ws.Range("J32").Value = 5045.5386
$ws.Range("K32").Value = 2498.0488
$ws.Range("L32").Value = 5045.5386
$ws.Range("M32").Value = -2211.0488
$ws.Range("N32").Value = -5619.5386
$ws.Range("H37").Value = 16300
$ws.Range("I37").Value = 13000
$ws.Range("J37").Value = 16960
$ws.Range("K37").Value = 13000
$ws.Range("L37").Value = 16960
$ws.Range("M37").Value = -12727
$ws.Range("N37").Value = -17506
$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30976
$ws.Range("H55").Value = 14000
$ws.Range("J55").Value = 14000
$ws.Range("L55").Value = 14000
$ws.Range("N55").Value = -14630
$ws.Range("H96").Value = 50000
$ws.Range("J96").Value = 50000
$ws.Range("L96").Value = 50000
$ws.Range("N96").Value = -55492
$ws.Range("H122").Value = 1001.86957
$ws.Range("I122").Value = 932.2941
$ws.Range("K122").Value = 2796.8823
$ws.Range("M122").Value = -346.8822999999998
$ws.Range("H132").Value = 2939.0605
$ws.Range("I132").Value = 2767.5356
$ws.Range("K132").Value = 8302.606800000001
$ws.Range("M132").Value = -5772.606800000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H94").Value = 1519.1852
$ws.Range("I94").Value = 1407.2
$ws.Range("J94").Value = 1659.1666
$ws.Range("K94").Value = 1407.2
$ws.Range("L94").Value = 1659.1666
$ws.Range("M94").Value = -956.2
$ws.Range("N94").Value = -2561.1666
$ws.Range("H105").Value = 2098.182
$ws.Range("I105").Value = 2081.9355
$ws.Range("K105").Value = 2081.9355
$ws.Range("M105").Value = -334.9355
$ws.Range("H134").Value = 2825
$ws.Range("I134").Value = 1650
$ws.Range("K134").Value = 4950
$ws.Range("M134").Value = -2415

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1033.2858
$ws.Range("I16").Value = 1008.25
$ws.Range("J16").Value = 1066.6666
$ws.Range("K16").Value = 1008.25
$ws.Range("L16").Value = 1066.6666
$ws.Range("M16").Value = -721.25
$ws.Range("N16").Value = -1640.6666
$ws.Range("H31").Value = 1204.1637
$ws.Range("I31").Value = 907.4211
$ws.Range("J31").Value = 1360.7778
$ws.Range("K31").Value = 907.4211
$ws.Range("L31").Value = 1360.7778
$ws.Range("M31").Value = -612.4211
$ws.Range("N31").Value = -1950.7778
$ws.Range("H34").Value = 1204.1637
$ws.Range("I34").Value = 907.4211
$ws.Range("J34").Value = 1360.7778
$ws.Range("K34").Value = 907.4211
$ws.Range("L34").Value = 1360.7778
$ws.Range("M34").Value = -705.4211
$ws.Range("N34").Value = -1764.7778
$ws.Range("H93").Value = 16399.6
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -53744
$ws.Range("H96").Value = 19159.666
$ws.Range("J96").Value = 19159.666
$ws.Range("L96").Value = 19159.666
$ws.Range("N96").Value = -24651.666
$ws.Range("H113").Value = 1033.2858
$ws.Range("I113").Value = 1008.25
$ws.Range("J113").Value = 1066.6666
$ws.Range("K113").Value = 1008.25
$ws.Range("L113").Value = 1066.6666
$ws.Range("M113").Value = 1161.75
$ws.Range("N113").Value = -5406.6666
$ws.Range("H132").Value = 3427.4
$ws.Range("I132").Value = 2866.4
$ws.Range("J132").Value = 5110.4
$ws.Range("K132").Value = 8599.200000000001
$ws.Range("L132").Value = 15331.2
$ws.Range("M132").Value = -6069.200000000001
$ws.Range("N132").Value = -20391.2

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 766.6667
$ws.Range("I32").Value = 300
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 900
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -617
$ws.Range("N32").Value = -3566
$ws.Range("H33").Value = 399
$ws.Range("I33").Value = 55
$ws.Range("J33").Value = 513.6667
$ws.Range("K33").Value = 330
$ws.Range("L33").Value = 3082.0002
$ws.Range("M33").Value = -47
$ws.Range("N33").Value = -3648.0002
$ws.Range("H98").Value = 469.27274
$ws.Range("J98").Value = 458.375
$ws.Range("L98").Value = 1375.125
$ws.Range("N98").Value = -4371.125
$ws.Range("H107").Value = 1489.9412
$ws.Range("J107").Value = 1735.2727
$ws.Range("L107").Value = 5205.8181
$ws.Range("N107").Value = -9045.8181
$ws.Range("H108").Value = 803.125
$ws.Range("I108").Value = 803.125
$ws.Range("K108").Value = 2409.375
$ws.Range("M108").Value = 470.625
$ws.Range("H139").Value = 11707
$ws.Range("J139").Value = 1000
$ws.Range("L139").Value = 3000
$ws.Range("N139").Value = -13280
$ws.Range("H140").Value = 2983.2222
$ws.Range("I140").Value = 653.94446
$ws.Range("J140").Value = 5312.5
$ws.Range("K140").Value = 1961.83338
$ws.Range("L140").Value = 15937.5
$ws.Range("M140").Value = 3218.16662
$ws.Range("N140").Value = -26297.5

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 19000
$ws.Range("J43").Value = 19000
$ws.Range("L43").Value = 19000
$ws.Range("N43").Value = -19302
$ws.Range("H46").Value = 19585.715
$ws.Range("J46").Value = 19585.715
$ws.Range("L46").Value = 19585.715
$ws.Range("N46").Value = -19897.715
$ws.Range("H70").Value = 4228.9165
$ws.Range("I70").Value = 4457.143
$ws.Range("J70").Value = 3909.4
$ws.Range("K70").Value = 4457.143
$ws.Range("L70").Value = 3909.4
$ws.Range("M70").Value = -4187.143
$ws.Range("N70").Value = -4449.4
$ws.Range("H73").Value = 4228.9165
$ws.Range("I73").Value = 4457.143
$ws.Range("J73").Value = 3909.4
$ws.Range("K73").Value = 4457.143
$ws.Range("L73").Value = 3909.4
$ws.Range("M73").Value = -3521.143
$ws.Range("N73").Value = -5781.4
$ws.Range("H92").Value = 12831.333
$ws.Range("J92").Value = 12831.333
$ws.Range("L92").Value = 12831.333
$ws.Range("N92").Value = -16575.333
$ws.Range("H126").Value = 1854641.2
$ws.Range("I126").Value = 2648209.2
$ws.Range("J126").Value = 2982.3333
$ws.Range("K126").Value = 7944627.600000001
$ws.Range("L126").Value = 8946.999899999999
$ws.Range("M126").Value = -7942157.600000001
$ws.Range("N126").Value = -13886.9999
$ws.Range("H132").Value = 1605325.6
$ws.Range("I132").Value = 2025990.4
$ws.Range("J132").Value = 6799.4
$ws.Range("K132").Value = 6077971.199999999
$ws.Range("L132").Value = 20398.2
$ws.Range("M132").Value = -6075441.199999999
$ws.Range("N132").Value = -25458.2
$ws.Range("H134").Value = 38874.25
$ws.Range("J134").Value = 38874.25
$ws.Range("L134").Value = 116622.75
$ws.Range("N134").Value = -121692.75
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 957.6842
$ws.Range("I46").Value = 684.2
$ws.Range("J46").Value = 1055.3572
$ws.Range("K46").Value = 684.2
$ws.Range("L46").Value = 1055.3572
$ws.Range("M46").Value = -496.2
$ws.Range("N46").Value = -1431.3572
$ws.Range("H60").Value = 20000
$ws.Range("J60").Value = 20000
$ws.Range("L60").Value = 20000
$ws.Range("N60").Value = -21018
$ws.Range("H132").Value = 3480
$ws.Range("I132").Value = 1499.8
$ws.Range("K132").Value = 4499.4
$ws.Range("M132").Value = -1969.4

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13332.667
$ws.Range("J41").Value = 13332.667
$ws.Range("L41").Value = 13332.667
$ws.Range("N41").Value = -14112.667
$ws.Range("H136").Value = 29243292
$ws.Range("I136").Value = 50508516
$ws.Range("K136").Value = 151525548
$ws.Range("M136").Value = -151522998
